$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text type on the Price/Volume columns so numeric-looking values
# (e.g. '1.001', '1.000') stay as text instead of being auto-converted to numbers.
$rng = $ws.Range("D2:E51")
$rng.NumberFormat = "@"

$ws.Range("D2").Value = "29.974.21"
$ws.Range("E2").Value = "  -0.13%  "

$ws.Range("D3").Value = "1.908.51"
$ws.Range("E3").Value = "  +0.36%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "0.7791"
$ws.Range("E5").Value = "  +4.71%  "

$ws.Range("D6").Value = "241.71"
$ws.Range("E6").Value = "  -0.17%  "

$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "0.3145"
$ws.Range("E8").Value = "  +2.41%  "

$ws.Range("D9").Value = "26.03"
$ws.Range("E9").Value = "  +1.65%  "

$ws.Range("D10").Value = "0.06873"
$ws.Range("E10").Value = "  -0.43%  "

$ws.Range("D11").Value = "0.07956"
$ws.Range("E11").Value = "  -1.12%  "

$ws.Range("D12").Value = "1.901.98"
$ws.Range("E12").Value = "  -0.10%  "

$ws.Range("D13").Value = "0.7390"
$ws.Range("E13").Value = "  -2.19%  "

$ws.Range("D14").Value = "5.188"
$ws.Range("E14").Value = "  -0.91%  "

$ws.Range("D15").Value = "92.67"
$ws.Range("E15").Value = "  +1.60%  "

$ws.Range("D16").Value = "29.973.69"
$ws.Range("E16").Value = "  -0.15%  "

$ws.Range("D17").Value = "13.88"
$ws.Range("E17").Value = "  -1.08%  "

$ws.Range("D18").Value = "5.853"
$ws.Range("E18").Value = "  -5.11%  "

$ws.Range("D19").Value = "244.92"
$ws.Range("E19").Value = "  +3.42%  "

$ws.Range("D20").Value = "0.000007720"
$ws.Range("E20").Value = "  -0.62%  "

$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  +0.01%  "

$ws.Range("D22").Value = "2.148.21"
$ws.Range("E22").Value = "  -0.27%  "

$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  -0.07%  "

$ws.Range("D24").Value = "6.841"
$ws.Range("E24").Value = "  -3.91%  "

$ws.Range("D25").Value = "168.52"
$ws.Range("E25").Value = "  +0.38%  "

$ws.Range("D26").Value = "9.250"
$ws.Range("E26").Value = "  -0.82%  "

$ws.Range("D27").Value = "0.1370"
$ws.Range("E27").Value = "  +7.62%  "

$ws.Range("D28").Value = "18.86"
$ws.Range("E28").Value = "  +0.32%  "

$ws.Range("D29").Value = "2.023"
$ws.Range("E29").Value = "  -1.35%  "

$ws.Range("E30").Value = "  +1.04%  "

$ws.Range("D31").Value = "1.515"
$ws.Range("E31").Value = "  -1.18%  "

$ws.Range("D32").Value = "4.303"
$ws.Range("E32").Value = "  +0.04%  "

$ws.Range("D33").Value = "4.074"
$ws.Range("E33").Value = "  +0.77%  "

$ws.Range("D34").Value = "0.05502"
$ws.Range("E34").Value = "  +4.22%  "

$ws.Range("D35").Value = "1.251"
$ws.Range("E35").Value = "  -2.49%  "

$ws.Range("D36").Value = "0.7310"
$ws.Range("E36").Value = "  -1.23%  "

$ws.Range("D37").Value = "2.729"
$ws.Range("E37").Value = "  +0.08%  "

$ws.Range("D38").Value = "0.01927"
$ws.Range("E38").Value = "  -0.92%  "

$ws.Range("D39").Value = "2.790"
$ws.Range("E39").Value = "  +0.93%  "

$ws.Range("D40").Value = "6.111"
$ws.Range("E40").Value = "  -2.24%  "

$ws.Range("D41").Value = "0.4405"
$ws.Range("E41").Value = "  -1.26%  "

$ws.Range("D42").Value = "71.73"
$ws.Range("E42").Value = "  -1.15%  "

$ws.Range("D43").Value = "1.001"
$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("D44").Value = "0.8387"
$ws.Range("E44").Value = "  +0.81%  "

$ws.Range("D45").Value = "1.868"
$ws.Range("E45").Value = "  -4.09%  "

$ws.Range("D46").Value = "100.31"
$ws.Range("E46").Value = "  -1.04%  "

$ws.Range("D47").Value = "7.500"
$ws.Range("E47").Value = "  -2.96%  "

$ws.Range("D48").Value = "9.725"
$ws.Range("E48").Value = "  -0.92%  "

$ws.Range("D49").Value = "978.87"
$ws.Range("E49").Value = "  +8.17%  "

$ws.Range("D50").Value = "2.057.00"
$ws.Range("E50").Value = "  -0.16%  "

$ws.Range("D51").Value = "36.09"
$ws.Range("E51").Value = "  -1.28%  "

# Restore the original (default) cell style now that the text values are set.
$rng.Style = "Normal"
